# daily auto push: 2026-01-24 22:33 UTC
# Insert a new daily record row right before the existing row 719,
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 719 (shifts rows 719.. down to 720..)
$ws.Rows.Item(719).Insert()

# Populate the newly inserted row with the new daily record.
# Force column A to text format so the date-like string is stored as
# plain text (matching the rest of the column) rather than being
# auto-converted into a date serial number.
$ws.Cells.Item(719, 1).NumberFormat = "@"
$ws.Cells.Item(719, 1).Value = "2026/01/25"
$ws.Cells.Item(719, 2).Value = "日"
$ws.Cells.Item(719, 3).Value = 5
$ws.Cells.Item(719, 4).Value = 168
